$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '{''chika'', ''naur'', ''shimenet'', ''anda'', ''andamhie'', ''eklabool''}'
$ws.Range("D3").Value = '{''chika'', ''λ'', ''naur'', ''shimenet'', ''anda'', ''andamhie'', ''eklabool''}'
$ws.Range("D4").Value = '{''chika'', ''λ'', ''naur'', ''shimenet'', ''anda'', ''andamhie'', ''eklabool''}'
$ws.Range("D5").Value = '{''λ'', ''naur''}'
$ws.Range("D6").Value = '{''chika'', ''anda'', ''andamhie'', ''eklabool'', ''shimenet''}'
$ws.Range("D7").Value = '{''andamhie'', ''eklabool'', ''chika'', ''anda''}'
$ws.Range("D8").Value = '{''['', ''('', ''='', ''λ''}'
$ws.Range("D9").Value = '{''['', ''(''}'
$ws.Range("D10").Value = '{''chika'', ''λ'', ''anda'', ''andamhie'', ''eklabool''}'
$ws.Range("D11").Value = '{''andamhie'', ''eklabool'', ''chika'', ''anda''}'
$ws.Range("D13").Value = '{''['', ''λ''}'
$ws.Range("D14").Value = '{''['', ''λ''}'
$ws.Range("D15").Value = '{''['', ''λ''}'
$ws.Range("D16").Value = '{''['', ''λ''}'
$ws.Range("D17").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D18").Value = '{''eme'', ''chika_literal'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D20").Value = '{''('', ''['', ''λ''}'
$ws.Range("D21").Value = '{''('', ''['', ''λ''}'
$ws.Range("D23").Value = '{''λ'', ''=''}'
$ws.Range("D25").Value = '{''λ'', ''=''}'
$ws.Range("D26").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D27").Value = '{''λ'', ''=''}'
$ws.Range("D28").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''{'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D30").Value = '{''chika'', ''anda'', ''λ'', ''andamhie'', ''eklabool'', ''shimenet''}'
$ws.Range("D31").Value = '{''push'', ''λ'', ''keri'', ''naur'', ''++'', ''--'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''andamhie'', ''id''}'
$ws.Range("D32").Value = '{''chika'', ''anda'', ''λ'', ''andamhie'', ''eklabool'', ''naur''}'
$ws.Range("D33").Value = '{''chika'', ''naur'', ''anda'', ''andamhie'', ''eklabool''}'
$ws.Range("D34").Value = '{''eme'', ''len'', ''-'', ''('', ''λ'', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D36").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D37").Value = '{''λ'', ''id''}'
$ws.Range("D38").Value = '{''push'', ''λ'', ''keri'', ''naur'', ''++'', ''--'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''andamhie'', ''id''}'
$ws.Range("D39").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D40").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D42").Value = '{''<='', ''//'', ''%'', ''>='', ''||'', ''-'', ''!='', ''λ'', ''*'', ''**'', ''<'', ''=='', ''>'', ''/'', ''&&'', ''+''}'
$ws.Range("D43").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D44").Value = '{''-'', ''λ'', ''!''}'
$ws.Range("D45").Value = '{''eme'', ''len'', ''('', ''++'', ''--'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D47").Value = '{''++'', ''--'', ''('', ''['', ''λ''}'
$ws.Range("D49").Value = '{''eme'', ''chika_literal'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D50").Value = '{''<='', ''//'', ''%'', ''>='', ''||'', ''-'', ''!='', ''*'', ''**'', ''<'', ''=='', ''>'', ''/'', ''&&'', ''+''}'
$ws.Range("D51").Value = '{''push'', ''λ'', ''keri'', ''naur'', ''++'', ''--'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''andamhie'', ''id''}'
$ws.Range("D52").Value = '{''amaccana'', ''push'', ''λ'', ''gogogo'', ''keri'', ''naur'', ''++'', ''--'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''andamhie'', ''id''}'
$ws.Range("D53").Value = '{''amaccana'', ''push'', ''λ'', ''gogogo'', ''keri'', ''naur'', ''++'', ''--'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''andamhie'', ''id''}'
$ws.Range("D54").Value = '{''+='', ''//='', ''*='', ''-='', ''['', ''**='', ''('', ''/='', ''%='', ''=''}'
$ws.Range("D55").Value = '{''+='', ''//='', ''*='', ''%='', ''-='', ''**='', ''='', ''/=''}'
$ws.Range("D56").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''{'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D57").Value = '{''chika'', ''anda'', ''andamhie'', ''id'', ''eklabool''}'
$ws.Range("D58").Value = '{''chika'', ''λ'', ''anda'', ''andamhie'', ''eklabool''}'
$ws.Range("D59").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D63").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D65").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D66").Value = '{''push'', ''λ'', ''keri'', ''naur'', ''++'', ''--'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''andamhie'', ''id''}'
$ws.Range("D68").Value = '{''amaccana'', ''push'', ''λ'', ''gogogo'', ''keri'', ''naur'', ''++'', ''--'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''andamhie'', ''id''}'
$ws.Range("D69").Value = '{''ganern'', ''λ''}'
$ws.Range("D70").Value = '{''ganern'', ''λ''}'
$ws.Range("D73").Value = '{''chika'', ''λ'', ''anda'', ''andamhie'', ''eklabool''}'
$ws.Range("D74").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D75").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D76").Value = '{''λ'', ''step''}'
$ws.Range("D77").Value = '{''eme'', ''len'', ''-'', ''('', ''λ'', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D78").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D79").Value = '{''amaccana'', ''push'', ''λ'', ''gogogo'', ''keri'', ''naur'', ''++'', ''--'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''andamhie'', ''id''}'
$ws.Range("D81").Value = '{''('', ''lang''}'
$ws.Range("D82").Value = '{''amaccana'', ''push'', ''λ'', ''gogogo'', ''keri'', ''naur'', ''++'', ''--'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''andamhie'', ''id''}'
$ws.Range("D84").Value = '{''λ'', ''betsung''}'
$ws.Range("D85").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D86").Value = '{''push'', ''λ'', ''keri'', ''naur'', ''++'', ''--'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''andamhie'', ''id''}'
$ws.Range("D88").Value = '{''λ'', ''betsung''}'
$ws.Range("D89").Value = '{''amaccana'', ''push'', ''λ'', ''gogogo'', ''keri'', ''naur'', ''++'', ''--'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''andamhie'', ''id''}'
$ws.Range("D91").Value = '{''λ'', ''ditech''}'
$ws.Range("D92").Value = '{''λ'', ''amaccana'', ''gogogo''}'
$ws.Range("D95").Value = '{''eme'', ''len'', ''-'', ''('', ''λ'', ''++'', ''--'', ''{'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
